# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = "70.889.29"; ForceText = $false }
    @{ Cell = "E2"; Value = "  -0.07%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "3.813.90"; ForceText = $false }
    @{ Cell = "E3"; Value = "  -0.79%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  -0.03%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "707.32"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +2.45%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "169.86"; ForceText = $true }
    @{ Cell = "E6"; Value = "  -1.63%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "3.799.98"; ForceText = $false }
    @{ Cell = "E7"; Value = "  -1.13%  "; ForceText = $false }
    @{ Cell = "E8"; Value = "  +0.08%  "; ForceText = $false }
    @{ Cell = "E9"; Value = "  -0.75%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.160"; ForceText = $true }
    @{ Cell = "E10"; Value = "  -1.77%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "7.36"; ForceText = $true }
    @{ Cell = "E11"; Value = "  -0.71%  "; ForceText = $false }
    @{ Cell = "E12"; Value = "  -1.01%  "; ForceText = $false }
    @{ Cell = "E13"; Value = "  -2.17%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "36.29"; ForceText = $true }
    @{ Cell = "E14"; Value = "  -0.67%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "4.463.41"; ForceText = $false }
    @{ Cell = "E15"; Value = "  -0.60%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "3.818.38"; ForceText = $false }
    @{ Cell = "E16"; Value = "  -1.17%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "70.914.53"; ForceText = $false }
    @{ Cell = "E17"; Value = "  -0.10%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "7.16"; ForceText = $true }
    @{ Cell = "E18"; Value = "  -0.66%  "; ForceText = $false }
    @{ Cell = "E19"; Value = "  +0.12%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "17.28"; ForceText = $true }
    @{ Cell = "E20"; Value = "  -2.68%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "493.32"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +0.96%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "10.57"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -4.63%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "0.728"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +1.12%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "85.59"; ForceText = $true }
    @{ Cell = "E24"; Value = "  +1.11%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "0.0000143"; ForceText = $true }
    @{ Cell = "E25"; Value = "  -1.76%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "10.53"; ForceText = $true }
    @{ Cell = "E26"; Value = "  -0.02%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "12.05"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -2.38%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "3.971.76"; ForceText = $false }
    @{ Cell = "E28"; Value = "  -0.64%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "0.999"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -0.20%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "2.06"; ForceText = $true }
    @{ Cell = "E30"; Value = "  -3.86%  "; ForceText = $false }
    @{ Cell = "E31"; Value = "  -0.93%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "7.37"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -3.27%  "; ForceText = $false }
    @{ Cell = "D33"; Value = "2.22"; ForceText = $true }
    @{ Cell = "E33"; Value = "  -3.65%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "29.18"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -1.75%  "; ForceText = $false }
    @{ Cell = "E35"; Value = "  -4.38%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "3.787.41"; ForceText = $false }
    @{ Cell = "E36"; Value = "  -0.21%  "; ForceText = $false }
    @{ Cell = "B37"; Value = "Binance-PegBSC-USD"; ForceText = $false }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; ForceText = $false }
    @{ Cell = "D37"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E37"; Value = "  +0.01%  "; ForceText = $false }
    @{ Cell = "B38"; Value = "Aptos"; ForceText = $false }
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; ForceText = $false }
    @{ Cell = "D38"; Value = "9.09"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -1.92%  "; ForceText = $false }
    @{ Cell = "D39"; Value = "0.101"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -2.01%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "1.04"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +3.01%  "; ForceText = $false }
    @{ Cell = "E41"; Value = "  -2.97%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "5.93"; ForceText = $true }
    @{ Cell = "E42"; Value = "  -2.06%  "; ForceText = $false }
    @{ Cell = "E43"; Value = "  -3.97%  "; ForceText = $false }
    @{ Cell = "E45"; Value = "  +0.04%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "163.92"; ForceText = $true }
    @{ Cell = "E46"; Value = "  -0.60%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "0.000309"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +1.10%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "425.56"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +2.82%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "48.76"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +0.23%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "8.71"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +0.22%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "0.294"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -2.68%  "; ForceText = $false }
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    if ($chg.ForceText) {
        # These values look like numbers (e.g. "707.32") but must be stored
        # as literal text to preserve exact formatting (trailing zeros, etc.)
        $rng.NumberFormat = "@"
    }
    $rng.Value = $chg.Value
}

Write-Output "Applied $($changes.Count) cell updates."
